$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D values from Excel's automatic numeric/date conversion
# by temporarily forcing Text format, then restoring the default style
# so the saved cell XML matches the original (no explicit style index).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.063.48'
$ws.Range('E2').Value = '  +3.07%  '
$ws.Range('D3').Value = '3.041.79'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '596.80'
$ws.Range('E5').Value = '  +1.55%  '
$ws.Range('D6').Value = '152.09'
$ws.Range('E6').Value = '  +6.98%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.037.30'
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('D10').Value = '6.37'
$ws.Range('E10').Value = '  +10.21%  '
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  +6.12%  '
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').Value = '0.0000235'
$ws.Range('E13').Value = '  +4.09%  '
$ws.Range('D14').Value = '34.99'
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = '3.551.73'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '63.078.33'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').Value = '7.04'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '3.049.60'
$ws.Range('E19').Value = '  +2.67%  '
$ws.Range('D20').Value = '455.20'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').Value = '14.26'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('D22').Value = '0.693'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').Value = '7.50'
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('D24').Value = '82.85'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  +4.86%  '
$ws.Range('D26').Value = '10.80'
$ws.Range('E26').Value = '  +8.62%  '
$ws.Range('D27').Value = '12.22'
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '2.72'
$ws.Range('E29').Value = '  +3.05%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.42'
$ws.Range('E30').Value = '  +8.21%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').Value = '2.16'
$ws.Range('E32').Value = '  +5.17%  '
$ws.Range('D33').Value = '27.69'
$ws.Range('E33').Value = '  +1.89%  '
$ws.Range('D34').Value = '0.110'
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('D35').Value = '0.0₃0853'
$ws.Range('E35').Value = '  +9.25%  '
$ws.Range('E36').Value = '  +2.49%  '
$ws.Range('D37').Value = '5.90'
$ws.Range('E37').Value = '  +2.70%  '
$ws.Range('D38').Value = '3.13'
$ws.Range('E38').Value = '  +13.14%  '
$ws.Range('D39').Value = '2.11'
$ws.Range('E39').Value = '  +2.53%  '
$ws.Range('D40').Value = '50.57'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '9.12'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('D43').Value = '0.294'
$ws.Range('E43').Value = '  +11.55%  '
$ws.Range('D44').Value = '41.03'
$ws.Range('E44').Value = '  +10.09%  '
$ws.Range('D45').Value = '394.06'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = '0.0357'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').Value = '2.749.03'
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('D48').Value = '133.09'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '2.21'
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.108'
$ws.Range('E51').Value = '  +0.69%  '

$ws.Range("D2:D51").Style = "Normal"
